# The edit inserts a new data row at row 56 in Sheet1 (pushing the former
# rows 56..131 down to 57..132) and populates the new row with the
# "Granada" market entry for Vega Modelo de Temuco described in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 56. Excel shifts rows 56-131
# down to 57-132, which is exactly what the diff shows happened to every
# row from 56 through 131 (each row's content became the content of the
# row below it).
$ws.Rows("56:56").Insert()

# Populate the newly inserted row 56 with the new record's values.
$ws.Range("A56").Value = 10
$ws.Range("B56").Value = "Vega Modelo de Temuco"
$ws.Range("C56").Value = "La Araucanía"
$ws.Range("D56").Value = 44740
$ws.Range("E56").Value = 9
$ws.Range("F56").Value = "Fruta"
$ws.Range("G56").Value = 100104
$ws.Range("H56").Value = "Frutos de pepita"
$ws.Range("I56").Value = 100104001
$ws.Range("J56").Value = "Granada"
$ws.Range("K56").Value = "Wonderfull"
$ws.Range("L56").Value = "Segunda"
$ws.Range("M56").Value = 200
$ws.Range("N56").Value = 13000
$ws.Range("O56").Value = 13000
$ws.Range("P56").Value = 13000
$ws.Range("Q56").Value = "$/bandeja 15 kilos granel"
$ws.Range("R56").Value = "Provincia de Limarí"
$ws.Range("S56").Value = 867
$ws.Range("T56").Value = 15
